# Auto-generated Excel COM-interop script applying the cryptos.xlsx diff
# (GitHub Actions "Updated cryptos list" commit): refreshes price (D) and
# hourly volume change % (E) columns for rows 2-51 of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.144.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.342.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.04%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.594"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.26%  "

$ws.Range("E9").Value = "  +5.34%  "

$ws.Range("E10").Value = "  +1.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "47.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.27%  "

$ws.Range("E12").Value = "  +2.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "699.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.884.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.11%  "

$ws.Range("E15").Value = "  +0.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "68.142.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.71%  "

$ws.Range("E17").Value = "  +0.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.337.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.33%  "

$ws.Range("E21").Value = "  +1.36%  "

$ws.Range("E22").Value = "  +1.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "101.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.41%  "

$ws.Range("E25").Value = "  +2.59%  "

$ws.Range("E26").Value = "  +1.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.52%  "

$ws.Range("E30").Value = "  -0.28%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "574.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.75%  "

$ws.Range("E32").Value = "  +1.26%  "

$ws.Range("E33").Value = "  +2.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.759.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.86%  "

$ws.Range("E35").Value = "  +0.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "56.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.20%  "

$ws.Range("E39").Value = "  +2.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.72%  "

$ws.Range("E41").Value = "  +0.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0681"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.336"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("E45").Value = "  +0.71%  "

$ws.Range("E46").Value = "  +2.28%  "

$ws.Range("E47").Value = "  +1.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.26%  "

$ws.Range("E49").Value = "  -0.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "130.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.48%  "
